# Updated cryptos list - applies latest price/volume scrape to the sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.041.20"
$ws.Range("E2").Value = "  -0.85%  "

$ws.Range("D3").Value = "2.550.56"
$ws.Range("E3").Value = "  +0.01%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.95"
$ws.Range("E5").Value = "  +2.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.99"
$ws.Range("E6").Value = "  -2.50%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("E8").Value = "  -0.50%  "

$ws.Range("E9").Value = "  -0.44%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.55"
$ws.Range("E10").Value = "  -3.11%  "

$ws.Range("E11").Value = "  -0.15%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.56"
$ws.Range("E13").Value = "  -3.01%  "

$ws.Range("D14").Value = "3.002.88"
$ws.Range("E14").Value = "  -0.20%  "

$ws.Range("D15").Value = "62.929.19"
$ws.Range("E15").Value = "  -0.89%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000144"
$ws.Range("E16").Value = "  -0.08%  "

$ws.Range("D17").Value = "2.542.28"
$ws.Range("E17").Value = "  -1.95%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.33"
$ws.Range("E18").Value = "  -2.88%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "338.52"
$ws.Range("E19").Value = "  -0.56%  "

$ws.Range("E20").Value = "  -1.18%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.76"
$ws.Range("E21").Value = "  -1.29%  "

$ws.Range("E22").Value = "  -0.19%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.61"
$ws.Range("E23").Value = "  -0.75%  "

$ws.Range("D24").Value = "2.678.99"
$ws.Range("E24").Value = "  +0.21%  "

$ws.Range("E25").Value = "  -0.48%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.61"
$ws.Range("E26").Value = "  +0.47%  "

$ws.Range("E27").Value = "  -2.61%  "

$ws.Range("E28").Value = "  -0.03%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.34"
$ws.Range("E29").Value = "  -2.31%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.70"
$ws.Range("E30").Value = "  +7.81%  "

$ws.Range("E31").Value = "  +3.10%  "

$ws.Range("D32").Value = "0.0₃0815"
$ws.Range("E32").Value = "  -1.45%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "177.87"
$ws.Range("E33").Value = "  -0.08%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "421.10"
$ws.Range("E34").Value = "  -0.42%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.55"
$ws.Range("E35").Value = "  -2.27%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.401"
$ws.Range("E36").Value = "  -1.38%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.09"
$ws.Range("E37").Value = "  -0.26%  "

$ws.Range("E38").Value = "  +0.03%  "

$ws.Range("E39").Value = "  -2.07%  "

$ws.Range("E40").Value = "  -2.25%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.14%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.76"
$ws.Range("E42").Value = "  -0.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "150.68"
$ws.Range("E43").Value = "  -2.27%  "

$ws.Range("E44").Value = "  -0.51%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.76"
$ws.Range("E45").Value = "  -1.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0538"
$ws.Range("E46").Value = "  +1.33%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.602"
$ws.Range("E47").Value = "  -1.69%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0968"
$ws.Range("E48").Value = "  -0.17%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0238"
$ws.Range("E49").Value = "  -0.45%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.28"
$ws.Range("E50").Value = "  -2.21%  "

$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.71"
$ws.Range("E51").Value = "  -5.90%  "

